$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 60: 2026/01/09, 逃离鸭科夫, 1140 (same style as the rows above it).
# A60 must stay plain text (not auto-parsed into a date serial), so force
# text format for the write, then clear the format back off again so the
# final style matches the existing "center/center" style (s="1") exactly.
$ws.Cells.Item(60, 1).NumberFormat = "@"
$ws.Cells.Item(60, 1).Value = "2026/01/09"
$ws.Cells.Item(60, 1).ClearFormats()

$ws.Cells.Item(60, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(60, 3).Value = 1140

$ws.Range("A60:C60").HorizontalAlignment = -4108
$ws.Range("A60:C60").VerticalAlignment = -4108
